# US351873 Add github parsing for SourceMinder_Product_Contact
#
# Inserts a new data row (GitHub / Workload Automation source) into the
# "Product Contact Matrix" sheet at row 679, pushing the existing rows
# 679-695 down to 680-696. The new row is highlighted with a green fill
# (to flag newly-parsed GitHub sources), matching the existing convention
# of colour-coding special rows (e.g. the light-blue fill already used on
# neighbouring rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Contact Matrix")

$insertRow = 679

# Shift rows 679:695 down by inserting a new blank row 679.
$ws.Rows.Item($insertRow).Insert()

$newRow = $ws.Rows.Item($insertRow)

$ws.Cells.Item($insertRow, 1).Value = "Workload Automation"
$ws.Cells.Item($insertRow, 2).Value = $null
$ws.Cells.Item($insertRow, 3).Value = "Active"
$ws.Cells.Item($insertRow, 4).Value = "GitHub"
$ws.Cells.Item($insertRow, 5).Value = "null"
$ws.Cells.Item($insertRow, 6).Value = "github-isl-01.ca.com/WLA"
$ws.Cells.Item($insertRow, 7).Value = '[{"PMFKEY":"bobsr01","TYPE":"Team","NAME":"Agents"},{"PMFKEY":"bobsr01","TYPE":"Repository","NAME":"ae"}]'
$ws.Cells.Item($insertRow, 8).Value = $null
$ws.Cells.Item($insertRow, 9).Value = $null
$ws.Cells.Item($insertRow, 10).Value = "bobsr01"

# Match row height / wrap formatting used throughout the sheet.
$newRow.RowHeight = 99.75
$rowRange = $ws.Range("A" + $insertRow + ":J" + $insertRow)
$rowRange.VerticalAlignment = -4160  # xlTop
$rowRange.WrapText = $true

# Highlight the new row with a green fill (new source-management-tool type).
# 0x00B050 -> OLE long (BGR) = B*65536 + G*256 + R
$rowRange.Interior.Color = 5287936

# Update the filter database defined name to include the extra row.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "='Product Contact Matrix'!`$A`$1:`$L`$696"

# Leave the cursor / scroll position where the edit happened, matching the
# author's saved view state.
$ws.Application.Goto($ws.Range("A" + $insertRow), $false)
